$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cells that previously held a 1 and must now be cleared (blank, matching
# the original cell's style which is retained by ClearContents).
$clearCells = @("Q2","P3","O4","N5","M6","L7","K8","J9","I10","H11","G12","F13","D15","O15","C16","P16","B17","Q17")
foreach ($addr in $clearCells) {
    $ws.Range($addr).ClearContents()
}

# Cells that are newly populated with the value 1.
$setCells = @("N10","I11","J11","L11","M11","N11","I12","J12","K12","M12","N12","H13","I13","J13","K13","L13","N13","C14","D14","F14","G14","H14","I14","J14","K14","L14","M14")
foreach ($addr in $setCells) {
    $ws.Range($addr).Value = 1
}
